# Database completed - pending usability testing and completion of documentation
#
# The underlying table data/order is unchanged; what changed is:
#  - punctuation clean-up in several cells (curly quotes -> straight quotes,
#    double spaces -> single spaces, en-dash -> hyphen) which also reshuffles
#    the shared-string table because Excel rewrites each edited cell's string
#    in-place (new unique strings get appended to the shared string table in
#    the order the cells were edited)
#  - a Data > Sort was (re-)applied on the table, which stamps a <sortState>
#  - the sheet's scroll/selection state changed (no more frozen/locked
#    topLeftCell, new selection at G15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix punctuation across the review text / one meal name -----------------
# Order matters: the shared-string table is rebuilt in the order cells are
# (re)written, so edit cells in the exact sequence that reproduces the
# target shared-string order.

$ws.Range("A1").Value = "ID"

$ws.Range("G3").Value = "This winter treat includes roast chicken, mixed 'stir fry veggies' and a delightful sauce. It is topped with crumpled up filo sheets and a generous sprinkling of sesame seeds. This 'pie' is extremely tasty and is best served piping hot. It is one of my mother's signature dishes, which means that if I give it anything less than a five-star rating, things will not end well."

$ws.Range("B3").Value = "Filo Chicken 'Pie'"

$ws.Range("G8").Value = "The trick to this one is to use homemade raisin bread (instead of shop-bought white / brown bread) and be generous with the amount of syrup / cinnamon and sugar added on the top. Provided the egg is properly cooked this is the ultimate in decadent breakfasts. It is a once a year treat that is totally worth every calorie."

$ws.Range("G2").Value = "This is a Sunday morning event and is one of my favourite foods. The omelet is made with one egg and filled with a generous serve of smoked salmon and cream cheese. In summer, the addition of fresh chives makes this meal extra-special. It is normally served with two slices of whole meal toast and butter."

$ws.Range("G5").Value = "To be fair I had this a *very* long time ago but still remember it as being too much rice and not nearly enough other fillings - so much so that I won't go there again."

$ws.Range("G4").Value = "For the uninitiated affogato is a scoop of vanilla ice cream 'drowned' in a shot of hot espresso. This combines two of my favourite things, ice cream and coffee. Provided the coffee (and ice cream) are of high quality this is the ultimate simple dessert. Tascas does this quite well and has earned a rating of four stars."

$ws.Range("G11").Value = "This is one of those 'make it in advance, freeze it and then eat it' jobs. It's a stew made of potato, carrots, lentils and spices. The good thing about it is that it is quick to reheat, healthy and filling. Having said that, it is not something you'd want to eat everyday, as it can be a bit boring."

$ws.Range("G7").Value = "The good thing about Subway is the choice when it comes to bread, salads and fillings. I'm somewhat ambivalent when it comes to the chicken sub - it's 'OK' but a six-inch is not enough, and a twelve inch is way too much. I like that the sandwich has less than six grams of fat, but I just can't get all that excited over what is basically a chicken sandwich."

$ws.Range("G6").Value = "This is one of my favourite 'eating out' dinners. The combination of flat bread, humus, lettuce, lamb and other random ingredients is really, really, good. The only slight ""problem"" is that the portion size is extremely generous which generally means having no room left for dessert."

$ws.Range("G9").Value = "This is probably the best sushi shop in Auckland. The prices are reasonable, the pieces are generously sized and the ratio of rice to filling is perfect (ie: loads of filling with a relatively small amount of rice). They have a surprisingly large range for the size of the shop and their salmon / tuna sushi is particularly good."

$ws.Range("G10").Value = "This one involves a large bowl of lettuce, avocado, salad dressing and smoked tinned tuna. It's healthy and surprisingly filling but if I'm honest, it's not all that inspiring. It tastes OK provided the lettuce is fresh, but it is nowhere near as nice as other options that are out there. The upside is that it is quick and easy to assemble."

# --- Re-apply the sort on the ID column (stamps a <sortState>) -------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2"))
$ws.Sort.SetRange($ws.Range("A1:G12"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- Update the view: scroll back to A1 and move the selection to G15 -----
$ws.Range("G15").Select()
